$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table "Condicion_Pacientes" currently spans A1:F43 (header + 42
# data rows). Add one more data row for 2020-05-12 (serial date 43963),
# growing the table (and the autoFilter / sheet dimension) to A1:F44.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$newRowIndex = 44

$ws.Cells.Item($newRowIndex, 1).Value = 43963
$ws.Cells.Item($newRowIndex, 2).Value = 305
$ws.Cells.Item($newRowIndex, 3).Value = 74
$ws.Cells.Item($newRowIndex, 4).Value = 237
$ws.Cells.Item($newRowIndex, 5).Value = 9
$ws.Cells.Item($newRowIndex, 6).Value = 17

# Match the date number format used by the rest of column A (reuse the
# existing style rather than creating a new numFmt entry).
$ws.Cells.Item($newRowIndex - 1, 1).Copy() | Out-Null
$ws.Cells.Item($newRowIndex, 1).PasteSpecial(-4122) | Out-Null

# Match the saved selection recorded for this workbook revision.
$ws.Range("D47").Select() | Out-Null
